# Fruta / hortaliza, semanal
# Insert a new weekly price record for Mango (Feria Lagunitas de Puerto Montt)
# at row 390, pushing the existing rows 390:471 down to 391:472.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 390 - this shifts rows 390:471
# down to 391:472 and copies formatting from the row above (keeps the
# date-formatted style on column D).
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(390, 1).Value  = 4
$ws.Cells.Item(390, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(390, 3).Value  = 'Los Lagos'
$ws.Cells.Item(390, 4).Value  = 45275
$ws.Cells.Item(390, 5).Value  = 10
$ws.Cells.Item(390, 6).Value  = 'Fruta'
$ws.Cells.Item(390, 7).Value  = 100108
$ws.Cells.Item(390, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(390, 9).Value  = 100108002
$ws.Cells.Item(390, 10).Value = 'Mango'
$ws.Cells.Item(390, 11).Value = 'Sin especificar'
$ws.Cells.Item(390, 12).Value = 'Primera'
$ws.Cells.Item(390, 13).Value = 200
$ws.Cells.Item(390, 14).Value = 14000
$ws.Cells.Item(390, 15).Value = 14000
$ws.Cells.Item(390, 16).Value = 14000
$ws.Cells.Item(390, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(390, 18).Value = 'Perú'
$ws.Cells.Item(390, 19).Value = 3500
$ws.Cells.Item(390, 20).Value = 4
